$d = $word.ActiveDocument

# 1) Table cell margins: left 128 -> 133 dxa (6.4pt -> 6.65pt) for both tables
foreach ($t in $d.Tables) {
    $t.LeftPadding = 6.65
}

# 2) Collapse "[recipient.postal_address;strconv=no]" -> "[attachmentRecipient.postal_address;strconv=no]"
$d.Content.Find.Execute("[recipient.postal_address;strconv=no]", $true, $false, $false, $false, $false,
                         $true, 1, $false, "[attachmentRecipient.postal_address;strconv=no]", 2)

# 3) Date text update (TIME field cached result)
$d.Content.Find.Execute("09/12/2019", $true, $false, $false, $false, $false,
                         $true, 1, $false, "02/01/2020", 2)

# 4) Collapse "[recipient.civility] [recipient.lastname]," -> "[attachmentRecipient.civility] [attachmentRecipient.lastname],"
$d.Content.Find.Execute("[recipient.civility] [recipient.lastname],", $true, $false, $false, $false, $false,
                         $true, 1, $false, "[attachmentRecipient.civility] [attachmentRecipient.lastname],", 2)

# 5) Collapse "Veuillez agréer, [recipient.civility], l'expression de nos salutations distinguées."
$d.Content.Find.Execute(
    "Veuillez agréer, [recipient.civility], l" + [char]0x2019 + "expression de nos salutations distingu" + [char]0x00E9 + "es.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Veuillez agréer, [attachmentRecipient.civility], l" + [char]0x2019 + "expression de nos salutations distingu" + [char]0x00E9 + "es.",
    2)

# 6) Header decorative line shape: nudge its size (wp:extent / a:xfrm a:ext)
$sec = $d.Sections(1)
$hdr = $sec.Headers(1)
for ($i = 1; $i -le $hdr.Shapes.Count; $i++) {
    $s = $hdr.Shapes($i)
    if ($s.Name -eq "Image1" -and $s.Height -lt 1) {
        $s.Width = 543.89996
        $s.Height = 0.35
    }
}

Write-Output "done"
